# Adds the new "Save" column (H) to the sheet:
#  - H1: header label "Save", formatted like the other header cells (B1:G1)
#  - H2: numeric value 1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell, text + same formatting as the existing header row.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# New data cell under it.
$ws.Range("H2").Value = 1
